{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that precedes it) that followed the\n// \"LOB1008: Ci\u00eancia, Tecnologia e Sociedade (Requisito fraco)\" line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst markerText = \"LOB1008: Ci\u00eancia, Tecnologia e Sociedade (Requisito fraco)\";\nconst footerTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === markerText) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  // The paragraph right after the marker is the blank separator paragraph,\n  // followed by the two footer paragraphs. Delete all three.\n  const toDelete = [];\n  if (items[markerIndex + 1] && items[markerIndex + 1].text === \"\") {\n    toDelete.push(items[markerIndex + 1]);\n  }\n  if (items[markerIndex + 2] && items[markerIndex + 2].text === footerTexts[0]) {\n    toDelete.push(items[markerIndex + 2]);\n  }\n  if (items[markerIndex + 3] && items[markerIndex + 3].text === footerTexts[1]) {\n    toDelete.push(items[markerIndex + 3]);\n  }\n\n  for (const para of toDelete) {\n    para.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that precedes it) that followed the\n# \"LOB1008: Ci\u00eancia, Tecnologia e Sociedade (Requisito fraco)\" line.\n\n$d = $word.ActiveDocument\n\n$markerText = \"LOB1008: Ci\u00eancia, Tecnologia e Sociedade (Requisito fraco)\"\n$footerText1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$footerText2 = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$markerIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*$markerText*\") {\n        $markerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -gt 0) {\n    $idxBlank = $markerIndex + 1\n    $idxFooter1 = $markerIndex + 2\n    $idxFooter2 = $markerIndex + 3\n\n    # Delete from the highest index down so earlier indices stay valid.\n    if ($idxFooter2 -le $d.Paragraphs.Count) {\n        $t2 = $d.Paragraphs.Item($idxFooter2).Range.Text\n        if ($t2 -like \"*$footerText2*\") {\n            $d.Paragraphs.Item($idxFooter2).Range.Delete()\n        }\n    }\n    if ($idxFooter1 -le $d.Paragraphs.Count) {\n        $t1 = $d.Paragraphs.Item($idxFooter1).Range.Text\n        if ($t1 -like \"*$footerText1*\") {\n            $d.Paragraphs.Item($idxFooter1).Range.Delete()\n        }\n    }\n    if ($idxBlank -le $d.Paragraphs.Count) {\n        $tb = $d.Paragraphs.Item($idxBlank).Range.Text\n        if ($tb.Trim() -eq \"\") {\n            $d.Paragraphs.Item($idxBlank).Range.Delete()\n        }\n    }\n}\n"}
